$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing June (06/2025) "Dia 2" total_venda value
$ws.Cells.Item(3, 2).Value = 27797.8

# Insert a new row for June (06/2025) "Dia 4" and shift the rows below it down
$ws.Rows.Item(5).Insert()
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 35372.46
$ws.Cells.Item(5, 3).Value = 6
$ws.Cells.Item(5, 4).Value = 2025
$ws.Cells.Item(5, 5).Value = "06/2025"
